# LOT2054.docx content update
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $result = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $result) {
        Write-Output ("NOT FOUND: " + $old)
    }
}

# 1) Activation date
Replace-Text "Ativação: 01/01/2018" "Ativação: 01/01/2025"

# 2) Objetivos (PT) - fix typos / wording
Replace-Text "Capacitar o aluno a identificar e gerenciar os riscos no ambiente de trabalho, com enfâse nos dos decorrentes das atividades em biotecnologia na industria e em laboratórios de pesquisa." "Capacitar o aluno a identificar e gerenciar os riscos no ambiente de trabalho, com ênfase nos decorrentes das atividades em biotecnologia na indústria e em laboratórios de pesquisa."

# 3) Objetivos (EN) - the paragraph exists but has an empty (italic) run; insert the English text
$objEnPara = $d.Paragraphs.Item(7)
if ($objEnPara.Range.Text.Trim().Length -eq 0) {
    $objEnPara.Range.InsertBefore("Empower the student to identify and manage risks in the workplace, with a focus on those arising from activities in biotechnology in industry and research laboratories.")
}

# 4) Programa resumido (PT) - trimmed sentence
Replace-Text "Introdução a segurança e medicina do trabalho; Conceitos técnico e aspectos legais em biossegurança; Risco em biossegurança ao trabalhador, a comunidade e o ambiente." "Introdução a segurança e medicina do trabalho; Conceitos técnico e aspectos legais em biossegurança;"

# 5) Programa resumido (EN)
Replace-Text "Introduction to occupational health and safety;Technical knowledge and legal aspects in biosafety; Biosafety risk to the worker, the community, and the environment;" "Introduction to occupational safety and medicine; Technical concepts and legal aspects in biosafety."

# 6) Programa (PT)
Replace-Text "Gestão em segurança do trabalho e estratégias de prevenção e controle dos riscos no ambiente de trabalho. Normas regulamentadoras;Avaliação de segurança biológica, química e radiológica em atividades de biotecnologia. Manejo e descarte de resíduo contaminado;Critérios e normas regulamentadoras para os diferentes níveis de biossegurança;Legislação para produção e manejo organismos geneticamente modificado e seus derivados;Biossegurança no manuseio de cobaias; Princípios de bioética;Estudos de casos" "Introdução a gestão em segurança do trabalho e estratégias de prevenção; Identificação e controle dos riscos ambientais (físicos, químicos e biológicos). Normas regulamentadoras;Classes de risco biológico, níveis de biossegurança e normas para a atividades de biotecnologia. Descarte e classificação de resíduo;Legislação para produção e manejo organismos geneticamente modificados (OGM) e seus derivados;Biossegurança no manuseio de cobaias; Princípios de bioética;Estudos de casos problemas e soluções"

# 7) Programa (EN)
Replace-Text "Work safety management and strategies for prevention and risks control in the workplace. Regulatory norms;Assessment of biological, chemical and radiological hazard in biotechnology procedure;Handling and disposal of contaminated waste;Biosafety level criteria and Biosecurity regulations; Regulation of genetically modified organism and its products;Biosafety in laboratory animal handling;Principles of Bioethics;Case studies;" "Introduction to occupational safety management and prevention strategies; Identification and control of occupational risks (physical, chemical, and biological). Regulatory standards; Biological risk classes, biosafety levels, and regulations for biotechnology activities. Waste disposal and classification; Legislation for the production and handling of genetically modified organisms (GMOs) and their derivatives; Biosafety in handling laboratory animals; Principles of bioethics; Case studies of problems and solutions"

# 8) Avaliação - Método
Replace-Text "Duas notas N1e N2 distribuídas ao longo do semestre. A composição das `"N`" fica critério dodocente." "A avaliação será composta por provas, exercícios, projetos, seminários, relatórios e estudos de casos que poderão compor as notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."

# 9) Avaliação - Critério
Replace-Text "MF = (N1 + N2)/2" "MF≥ 5,0 para aprovação 5,0. Prova de recuperação para alunos com 3,0≤MF<5,0"

# 10) Avaliação - Norma de recuperação
Replace-Text "NF = (MF + PR)/2, onde PR é uma prova de recuperação.Prova de Recuperação (PR) para alunos com Média Final (MF) maior ou igual a 3,0 e menor doque 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final (NF) igual ou maior do que 5,0." "(MF+RC)/2 ≥ 5,0 para aprovação, onde RC é uma prova de recuperação a ser aplicada"

# 11) Bibliografia
Replace-Text "1- ABIQUIM . Código de saúde e segurança do trabalhadorguia de implantação.ABIQUIM, 1994. 2- Gonçalves, L. B. Gestão de segurança e medicina do trabalho normas regulamentadoras e fator acidentário de prevenção. Cenofisco, 2011.3- Paoleschi, B. Guia Prático de Segurança do Trabalho. Érica, 2009.4- Bisnfeld, P. C. Biossegurança em Biotecnologia. Interciência, 2004.5- Teixeira. P. e Valle, E. Biossegurança: uma abordagem multidisciplinar, 2002.6- Biosafety in Microbiological and Biomedical Laboratories, 5 ed. U.S. Health Department, 2013.8- Biosecurity, 1ed, 2013." "1-Binsfeld, P. C. Fundamentos Técnicos e o Sistema Nacional de Biossegurança em Biotecnologia. Interciência, 1ª edição 2015.2-Gonçalves Simão, L. B. Gestão de Segurança e Medicina do Trabalho, Normas Regulamentadoras e Fator Acidentário de Prevenção. Cenofisco, 1ª edição 2015.3-Hirata, M.H., Mancini Filho, J. Hirata, R. D. C. Manual de biossegurança.  Editora Manole. 3ª edição 2016.4- Semplici, S. Onze Teses de Bioética. Editora Ideias e Letras;1ª edição 2014"

Write-Output "done"
